# "fix all face-level issues of M2"
#
# Several text runs in the deck were generated with literal embedded
# newlines (hard line breaks baked into the run text) instead of letting
# PowerPoint wrap the paragraph naturally. This collapses every embedded
# newline in those runs into a single space so the text flows as one
# paragraph. Shapes use <a:spAutoFit/>, so touching their text makes the
# host re-layout/re-measure them (height changes are a side effect of the
# text edit, not something we set directly).

$p = $ppt.ActivePresentation

function Fix-Newlines {
    param($shape)
    $tr = $shape.TextFrame.TextRange
    $old = $tr.Text
    $new = $old.Replace("`r`n", "`n").Replace("`n", " ")
    if ($new -ne $old) {
        $tr.Text = $new
    }
}

# Slide 5 - "Code Block Support": Python + JavaScript code samples
$slide5 = $p.Slides.Item(5)
Fix-Newlines $slide5.Shapes.Item("TextBox 3")   # def fibonacci(n): ...
Fix-Newlines $slide5.Shapes.Item("TextBox 5")   # async function fetchUserData ...

# Slide 6 - SQL code sample (only shape on the slide)
$slide6 = $p.Slides.Item(6)
Fix-Newlines $slide6.Shapes.Item("TextBox 1")   # -- Complex query with joins ...

# Slide 11 - "Recent Improvements": Problem/Solution + Before/After call-outs
$slide11 = $p.Slides.Item(11)
Fix-Newlines $slide11.Shapes.Item("TextBox 3")  # Problem: Columns distributed ... / Solution: ...
Fix-Newlines $slide11.Shapes.Item("TextBox 4")  # Before: ... due to equal distribution / After: ...
Fix-Newlines $slide11.Shapes.Item("TextBox 6")  # Problem: Black borders invisible ... / Solution: ...
Fix-Newlines $slide11.Shapes.Item("TextBox 7")  # Dark Theme: ... for visibility / Default Theme: ...

# Slide 14 - "Technical Details": Python API usage sample
$slide14 = $p.Slides.Item(14)
Fix-Newlines $slide14.Shapes.Item("TextBox 4")  # from slide_generator.generator import SlideGenerator ...
